$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated / new values for rows 2-11 (Q0-Q9)
$data = @(
    @{ row = 2;  label = "Q0"; B = 0.2156134072807172;  C = 0.5104705812318155; D = 0.5597049646376884; E = 0.7481343225903276; F = 0.7243069924489576; G = 46 },
    @{ row = 3;  label = "Q1"; B = 0.07078259347676455; C = 0.3937967543584671; D = 0.2924574528699269; E = 0.5407933550534131; F = 0.5421993821439078; G = 45 },
    @{ row = 4;  label = "Q2"; B = 0.2022458033653479;  C = 0.5029758090258937; D = 0.6072955914828451; E = 0.7792917242489138; F = 0.761291099874567;  G = 44 },
    @{ row = 5;  label = "Q3"; B = 0.1056269426232762;  C = 0.3853209831470272; D = 0.2947523084012032; E = 0.5429109580780288; F = 0.5388390533604958; G = 43 },
    @{ row = 6;  label = "Q4"; B = 0.2315116146747508;  C = 0.5190355693975161; D = 0.6636775128766533; E = 0.8146640490881215; F = 0.7905440420027581; G = 42 },
    @{ row = 7;  label = "Q5"; B = 0.09265415981041521; C = 0.3465857693131973; D = 0.2074256174716431; E = 0.4554400262072308; F = 0.4514552522066898; G = 41 },
    @{ row = 8;  label = "Q6"; B = 0.2923186457898964;  C = 0.5389143003675101; D = 0.734479721597224;   E = 0.8570179237316008; F = 0.8158868265665942; G = 40 },
    @{ row = 9;  label = "Q7"; B = 0.1519280332023545;  C = 0.335957238970132;   D = 0.2134926424701896; E = 0.4620526403670794; F = 0.4420648348703681; G = 39 },
    @{ row = 10; label = "Q8"; B = 0.2798477600321015;  C = 0.4917781050680664; D = 0.5928569493473916; E = 0.7699720445233006; F = 0.7269447888063731; G = 38 },
    @{ row = 11; label = "Q9"; B = 0.1565110834604237;  C = 0.3372223941819412; D = 0.198541515633084;  E = 0.4455799766967586; F = 0.4229425514680761; G = 37 }
)

# Copy the format of the existing Q0 label cell (A2) so the new row labels
# (A6:A11) pick up the same bold / centered / bordered style.
$labelFormatSource = $ws.Cells.Item(2, 1)
$labelFormatSource.Copy()

foreach ($entry in $data) {
    $r = $entry.row
    $ws.Cells.Item($r, 1).Value = $entry.label
    if ($r -gt 5) {
        $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    }
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}
